$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------
# Step 1: Swap match-data columns (F:V) between specific row pairs.
# The identifying columns (A index, B/C/D/E metadata) stay on their original
# row; only the match details (teams, scores, odds, timestamps, url) swap.
# ----------------------------------------------------------------------------

    # Swap row 66 and 67 (columns F:V)
    $ws.Range("F66").Value = "Instituto"
    $ws.Range("F67").Value = "Racing Club"
    $ws.Range("G66").Value = 3
    $ws.Range("G67").Value = 2
    $ws.Range("H66").Value = "Colon Santa Fe"
    $ws.Range("H67").Value = "Newells Old Boys"
    $ws.Range("I66").Value = 1
    $ws.Range("I67").Value = 1
    $ws.Range("J66").Value = 2.04
    $ws.Range("J67").Value = 2.62
    $ws.Range("K66").Value = "17/09/2023 01:12"
    $ws.Range("K67").Value = "16/09/2023 21:12"
    $ws.Range("L66").Value = 2.18
    $ws.Range("L67").Value = 2.35
    $ws.Range("M66").Value = "21/09/2023 01:56"
    $ws.Range("M67").Value = "21/09/2023 01:55"
    $ws.Range("N66").Value = 3.19
    $ws.Range("N67").Value = 2.94
    $ws.Range("O66").Value = "17/09/2023 01:12"
    $ws.Range("O67").Value = "16/09/2023 21:12"
    $ws.Range("P66").Value = 3
    $ws.Range("P67").Value = 3.13
    $ws.Range("Q66").Value = "21/09/2023 01:56"
    $ws.Range("Q67").Value = "21/09/2023 01:55"
    $ws.Range("R66").Value = 4.18
    $ws.Range("R67").Value = 3.12
    $ws.Range("S66").Value = "17/09/2023 01:12"
    $ws.Range("S67").Value = "16/09/2023 21:12"
    $ws.Range("T66").Value = 4.13
    $ws.Range("T67").Value = 3.47
    $ws.Range("U66").Value = "21/09/2023 01:56"
    $ws.Range("U67").Value = "21/09/2023 01:55"
    $ws.Range("V66").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/instituto-colon-santa-fe/A90iIGQi/"
    $ws.Range("V67").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/racing-club-newells-old-boys/4Wym7jIj/"

    # Swap row 82 and 83 (columns F:V)
    $ws.Range("F82").Value = "Godoy Cruz"
    $ws.Range("F83").Value = "Colon Santa Fe"
    $ws.Range("G82").Value = 1
    $ws.Range("G83").Value = 3
    $ws.Range("H82").Value = "Racing Club"
    $ws.Range("H83").Value = "Argentinos Jrs"
    $ws.Range("I82").Value = 1
    $ws.Range("I83").Value = 1
    $ws.Range("J82").Value = 2.34
    $ws.Range("J83").Value = 3.03
    $ws.Range("K82").Value = "24/09/2023 05:42"
    $ws.Range("K83").Value = "21/09/2023 23:42"
    $ws.Range("L82").Value = 2.29
    $ws.Range("L83").Value = 2.95
    $ws.Range("M82").Value = "25/09/2023 23:29"
    $ws.Range("M83").Value = "25/09/2023 23:25"
    $ws.Range("N82").Value = 3.2
    $ws.Range("N83").Value = 3.09
    $ws.Range("O82").Value = "24/09/2023 05:42"
    $ws.Range("O83").Value = "21/09/2023 23:42"
    $ws.Range("P82").Value = 3.36
    $ws.Range("P83").Value = 3.16
    $ws.Range("Q82").Value = "25/09/2023 23:29"
    $ws.Range("Q83").Value = "25/09/2023 23:20"
    $ws.Range("R82").Value = 3.3
    $ws.Range("R83").Value = 2.46
    $ws.Range("S82").Value = "24/09/2023 05:42"
    $ws.Range("S83").Value = "21/09/2023 23:42"
    $ws.Range("T82").Value = 3.35
    $ws.Range("T83").Value = 2.64
    $ws.Range("U82").Value = "25/09/2023 23:29"
    $ws.Range("U83").Value = "25/09/2023 23:29"
    $ws.Range("V82").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/godoy-cruz-racing-club/8OvDrDvc/"
    $ws.Range("V83").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/colon-santa-fe-argentinos-jrs/0IHoSifG/"

    # Swap row 96 and 97 (columns F:V)
    $ws.Range("F96").Value = "Godoy Cruz"
    $ws.Range("F97").Value = "Barracas Central"
    $ws.Range("G96").Value = 1
    $ws.Range("G97").Value = 1
    $ws.Range("H96").Value = "Instituto"
    $ws.Range("H97").Value = "Sarmiento Junin"
    $ws.Range("I96").Value = 1
    $ws.Range("I97").Value = 1
    $ws.Range("J96").Value = 2.25
    $ws.Range("J97").Value = 2.3
    $ws.Range("K96").Value = "25/09/2023 22:42"
    $ws.Range("K97").Value = "26/09/2023 01:12"
    $ws.Range("L96").Value = 2.49
    $ws.Range("L97").Value = 2.14
    $ws.Range("M96").Value = "02/10/2023 23:29"
    $ws.Range("M97").Value = "02/10/2023 23:29"
    $ws.Range("N96").Value = 3.04
    $ws.Range("N97").Value = 3.05
    $ws.Range("O96").Value = "25/09/2023 22:42"
    $ws.Range("O97").Value = "26/09/2023 01:12"
    $ws.Range("P96").Value = 2.86
    $ws.Range("P97").Value = 3.03
    $ws.Range("Q96").Value = "02/10/2023 23:29"
    $ws.Range("Q97").Value = "02/10/2023 23:29"
    $ws.Range("R96").Value = 3.7
    $ws.Range("R97").Value = 3.54
    $ws.Range("S96").Value = "25/09/2023 22:42"
    $ws.Range("S97").Value = "26/09/2023 01:12"
    $ws.Range("T96").Value = 3.54
    $ws.Range("T97").Value = 4.22
    $ws.Range("U96").Value = "02/10/2023 23:29"
    $ws.Range("U97").Value = "02/10/2023 23:29"
    $ws.Range("V96").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/godoy-cruz-instituto/t4LKO1te/"
    $ws.Range("V97").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/barracas-central-sarmiento-junin/WA73IQJL/"

    # Swap row 98 and 99 (columns F:V)
    $ws.Range("F98").Value = "Platense"
    $ws.Range("F99").Value = "Atl. Tucuman"
    $ws.Range("G98").Value = 0
    $ws.Range("G99").Value = 0
    $ws.Range("H98").Value = "Argentinos Jrs"
    $ws.Range("H99").Value = "Central Cordoba"
    $ws.Range("I98").Value = 0
    $ws.Range("I99").Value = 0
    $ws.Range("J98").Value = 3.31
    $ws.Range("J99").Value = 1.93
    $ws.Range("K98").Value = "26/09/2023 01:12"
    $ws.Range("K99").Value = "26/09/2023 01:12"
    $ws.Range("L98").Value = 3.28
    $ws.Range("L99").Value = 1.91
    $ws.Range("M98").Value = "03/10/2023 01:51"
    $ws.Range("M99").Value = "03/10/2023 01:51"
    $ws.Range("N98").Value = 2.97
    $ws.Range("N99").Value = 3.29
    $ws.Range("O98").Value = "26/09/2023 01:12"
    $ws.Range("O99").Value = "26/09/2023 01:12"
    $ws.Range("P98").Value = 2.93
    $ws.Range("P99").Value = 3.2
    $ws.Range("Q98").Value = "03/10/2023 01:51"
    $ws.Range("Q99").Value = "03/10/2023 01:51"
    $ws.Range("R98").Value = 2.47
    $ws.Range("R99").Value = 4.52
    $ws.Range("S98").Value = "26/09/2023 01:12"
    $ws.Range("S99").Value = "26/09/2023 01:12"
    $ws.Range("T98").Value = 2.59
    $ws.Range("T99").Value = 5.03
    $ws.Range("U98").Value = "03/10/2023 01:51"
    $ws.Range("U99").Value = "03/10/2023 01:51"
    $ws.Range("V98").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/platense-argentinos-jrs/CrrxXPYE/"
    $ws.Range("V99").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/atl-tucuman-central-cordoba-santiago-del-estero/2DMGPsRl/"

    # Swap row 102 and 103 (columns F:V)
    $ws.Range("F102").Value = "Estudiantes L.P."
    $ws.Range("F103").Value = "Lanus"
    $ws.Range("G102").Value = 0
    $ws.Range("G103").Value = 0
    $ws.Range("H102").Value = "Godoy Cruz"
    $ws.Range("H103").Value = "Defensa y Justicia"
    $ws.Range("I102").Value = 1
    $ws.Range("I103").Value = 2
    $ws.Range("J102").Value = 1.86
    $ws.Range("J103").Value = 1.81
    $ws.Range("K102").Value = "02/10/2023 22:42"
    $ws.Range("K103").Value = "01/10/2023 02:42"
    $ws.Range("L102").Value = 2.17
    $ws.Range("L103").Value = 1.9
    $ws.Range("M102").Value = "07/10/2023 23:59"
    $ws.Range("M103").Value = "07/10/2023 23:58"
    $ws.Range("N102").Value = 3.37
    $ws.Range("N103").Value = 3.58
    $ws.Range("O102").Value = "02/10/2023 22:42"
    $ws.Range("O103").Value = "01/10/2023 02:42"
    $ws.Range("P102").Value = 3.06
    $ws.Range("P103").Value = 3.47
    $ws.Range("Q102").Value = "07/10/2023 23:59"
    $ws.Range("Q103").Value = "07/10/2023 23:58"
    $ws.Range("R102").Value = 4.33
    $ws.Range("R103").Value = 4.7
    $ws.Range("S102").Value = "02/10/2023 22:42"
    $ws.Range("S103").Value = "01/10/2023 02:42"
    $ws.Range("T102").Value = 4.04
    $ws.Range("T103").Value = 4.55
    $ws.Range("U102").Value = "07/10/2023 23:59"
    $ws.Range("U103").Value = "07/10/2023 23:58"
    $ws.Range("V102").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/estudiantes-l-p-godoy-cruz/x6gob2XD/"
    $ws.Range("V103").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/lanus-defensa-y-justicia/EX7jxLAf/"

    # Swap row 110 and 111 (columns F:V)
    $ws.Range("F110").Value = "Barracas Central"
    $ws.Range("F111").Value = "Arsenal Sarandi"
    $ws.Range("G110").Value = 2
    $ws.Range("G111").Value = 0
    $ws.Range("H110").Value = "Colon Santa Fe"
    $ws.Range("H111").Value = "Banfield"
    $ws.Range("I110").Value = 1
    $ws.Range("I111").Value = 0
    $ws.Range("J110").Value = 2.5
    $ws.Range("J111").Value = 3.26
    $ws.Range("K110").Value = "02/10/2023 22:42"
    $ws.Range("K111").Value = "02/10/2023 22:12"
    $ws.Range("L110").Value = 2.6
    $ws.Range("L111").Value = 3.55
    $ws.Range("M110").Value = "09/10/2023 22:36"
    $ws.Range("M111").Value = "09/10/2023 22:59"
    $ws.Range("N110").Value = 3.13
    $ws.Range("N111").Value = 2.96
    $ws.Range("O110").Value = "02/10/2023 22:42"
    $ws.Range("O111").Value = "02/10/2023 22:12"
    $ws.Range("P110").Value = 2.84
    $ws.Range("P111").Value = 2.9
    $ws.Range("Q110").Value = "09/10/2023 22:36"
    $ws.Range("Q111").Value = "09/10/2023 22:58"
    $ws.Range("R110").Value = 2.93
    $ws.Range("R111").Value = 2.51
    $ws.Range("S110").Value = "02/10/2023 22:42"
    $ws.Range("S111").Value = "02/10/2023 22:12"
    $ws.Range("T110").Value = 2.91
    $ws.Range("T111").Value = 2.46
    $ws.Range("U110").Value = "09/10/2023 22:36"
    $ws.Range("U111").Value = "09/10/2023 22:53"
    $ws.Range("V110").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/barracas-central-colon-santa-fe/4x7Fpq3E/"
    $ws.Range("V111").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/arsenal-sarandi-banfield/lUOSMuB7/"


# ----------------------------------------------------------------------------
# Step 2: Append three new match rows (162-164) at the end of the sheet.
# First copy the formatting (styles) of the last existing row (161) down to
# the new rows, then overwrite the values for each new row.
# ----------------------------------------------------------------------------

$ws.Range("A161:V161").Copy($ws.Range("A162:V162"))
$ws.Range("A161:V161").Copy($ws.Range("A163:V163"))
$ws.Range("A161:V161").Copy($ws.Range("A164:V164"))

    # New row 162
    $ws.Range("A162").Value = 161
    $ws.Range("B162").Value = "argentina"
    $ws.Range("C162").Value = "copa-de-la-liga-profesional"
    $ws.Range("E162").Value = 45235.9375
    $ws.Range("F162").Value = "Racing Club"
    $ws.Range("G162").Value = 1
    $ws.Range("H162").Value = "Central Cordoba"
    $ws.Range("I162").Value = 1
    $ws.Range("J162").Value = 1.67
    $ws.Range("K162").Value = "02/11/2023 00:12"
    $ws.Range("L162").Value = 1.6
    $ws.Range("M162").Value = "05/11/2023 22:26"
    $ws.Range("N162").Value = 3.68
    $ws.Range("O162").Value = "02/11/2023 00:12"
    $ws.Range("P162").Value = 3.96
    $ws.Range("Q162").Value = "05/11/2023 22:26"
    $ws.Range("R162").Value = 5.15
    $ws.Range("S162").Value = "02/11/2023 00:12"
    $ws.Range("T162").Value = 6.28
    $ws.Range("U162").Value = "05/11/2023 22:26"
    $ws.Range("V162").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/racing-club-central-cordoba-santiago-del-estero/UsaAPySs/"

    # New row 163
    $ws.Range("A163").Value = 162
    $ws.Range("B163").Value = "argentina"
    $ws.Range("C163").Value = "copa-de-la-liga-profesional"
    $ws.Range("E163").Value = 45236.02083333334
    $ws.Range("F163").Value = "Estudiantes L.P."
    $ws.Range("G163").Value = 2
    $ws.Range("H163").Value = "Defensa y Justicia"
    $ws.Range("I163").Value = 1
    $ws.Range("J163").Value = 1.88
    $ws.Range("K163").Value = "02/11/2023 00:12"
    $ws.Range("L163").Value = 1.85
    $ws.Range("M163").Value = "06/11/2023 00:29"
    $ws.Range("N163").Value = 3.39
    $ws.Range("O163").Value = "02/11/2023 00:12"
    $ws.Range("P163").Value = 3.44
    $ws.Range("Q163").Value = "06/11/2023 00:29"
    $ws.Range("R163").Value = 4.25
    $ws.Range("S163").Value = "02/11/2023 00:12"
    $ws.Range("T163").Value = 4.85
    $ws.Range("U163").Value = "06/11/2023 00:29"
    $ws.Range("V163").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/estudiantes-l-p-defensa-y-justicia/nJgFOesl/"

    # New row 164
    $ws.Range("A164").Value = 163
    $ws.Range("B164").Value = "argentina"
    $ws.Range("C164").Value = "copa-de-la-liga-profesional"
    $ws.Range("E164").Value = 45236.04166666666
    $ws.Range("F164").Value = "Newells Old Boys"
    $ws.Range("G164").Value = 0
    $ws.Range("H164").Value = "Sarmiento Junin"
    $ws.Range("I164").Value = 1
    $ws.Range("J164").Value = 1.68
    $ws.Range("K164").Value = "31/10/2023 01:42"
    $ws.Range("L164").Value = 1.74
    $ws.Range("M164").Value = "06/11/2023 00:59"
    $ws.Range("N164").Value = 3.54
    $ws.Range("O164").Value = "31/10/2023 01:42"
    $ws.Range("P164").Value = 3.43
    $ws.Range("Q164").Value = "06/11/2023 00:59"
    $ws.Range("R164").Value = 5.83
    $ws.Range("S164").Value = "31/10/2023 01:42"
    $ws.Range("T164").Value = 5.96
    $ws.Range("U164").Value = "06/11/2023 00:59"
    $ws.Range("V164").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/newells-old-boys-sarmiento-junin/KWiNMZB0/"

# Note: the sheet's <dimension> (used range) is recomputed automatically by
# Excel when the workbook is saved, so it will correctly become A1:V164.
